# This script updates the '想去人数' (want-to-go count) column (F) values
# across all four sheets, reflecting refreshed scrape counts.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # 展览
$ws.Range("F2").Value = 254
$ws.Range("F5").Value = 362
$ws.Range("F6").Value = 558
$ws.Range("F7").Value = 47
$ws.Range("F9").Value = 273
$ws.Range("F10").Value = 386
$ws.Range("F12").Value = 694
$ws.Range("F13").Value = 761
$ws.Range("F14").Value = 1517
$ws.Range("F15").Value = 1517
$ws.Range("F16").Value = 891
$ws.Range("F17").Value = 29
$ws.Range("F18").Value = 1355
$ws.Range("F19").Value = 164
$ws.Range("F20").Value = 326
$ws.Range("F23").Value = 103
$ws.Range("F24").Value = 6629
$ws.Range("F25").Value = 4990
$ws.Range("F27").Value = 487
$ws.Range("F28").Value = 207
$ws.Range("F29").Value = 202
$ws.Range("F32").Value = 1285
$ws.Range("F33").Value = 194
$ws.Range("F35").Value = 614
$ws.Range("F38").Value = 249
$ws.Range("F40").Value = 148
$ws.Range("F41").Value = 62

$ws = $wb.Worksheets.Item(2)  # 演出
$ws.Range("F10").Value = 12
$ws.Range("F15").Value = 51
$ws.Range("F18").Value = 241

$ws = $wb.Worksheets.Item(3)  # 本地生活
$ws.Range("F3").Value = 2460
$ws.Range("F5").Value = 59

$ws = $wb.Worksheets.Item(4)  # 全部类型
$ws.Range("F2").Value = 254
$ws.Range("F7").Value = 59
$ws.Range("F8").Value = 362
$ws.Range("F9").Value = 558
$ws.Range("F10").Value = 47
$ws.Range("F12").Value = 273
$ws.Range("F14").Value = 386
$ws.Range("F16").Value = 694
$ws.Range("F17").Value = 761
$ws.Range("F18").Value = 1517
$ws.Range("F19").Value = 1517
$ws.Range("F20").Value = 891
$ws.Range("F21").Value = 29
$ws.Range("F22").Value = 1355
$ws.Range("F23").Value = 164
$ws.Range("F24").Value = 326
$ws.Range("F26").Value = 103
$ws.Range("F29").Value = 6629
$ws.Range("F30").Value = 4990
$ws.Range("F32").Value = 202
$ws.Range("F34").Value = 1285
$ws.Range("F35").Value = 194
$ws.Range("F37").Value = 12
$ws.Range("F39").Value = 614
$ws.Range("F41").Value = 51
$ws.Range("F43").Value = 249
$ws.Range("F44").Value = 148
$ws.Range("F45").Value = 62
$ws.Range("F49").Value = 241

